{"js": "// Update the date label and all three-digit x one-digit multiplication\n// expressions in the practice sheet table to the new day's values.\nconst replacements = [\n  [\"2025-07-10 Thursday\", \"2025-07-11 Friday\"],\n  [\"122\u00d72=\", \"305\u00d73=\"],\n  [\"868\u00d79=\", \"470\u00d79=\"],\n  [\"327\u00d78=\", \"957\u00d77=\"],\n  [\"200\u00d72=\", \"296\u00d73=\"],\n  [\"116\u00d77=\", \"770\u00d78=\"],\n  [\"685\u00d79=\", \"324\u00d79=\"],\n  [\"891\u00d73=\", \"628\u00d75=\"],\n  [\"603\u00d77=\", \"231\u00d76=\"],\n  [\"964\u00d75=\", \"814\u00d78=\"],\n  [\"322\u00d78=\", \"780\u00d75=\"],\n  [\"681\u00d75=\", \"310\u00d72=\"],\n  [\"587\u00d78=\", \"817\u00d79=\"],\n  [\"325\u00d78=\", \"180\u00d74=\"],\n  [\"602\u00d77=\", \"965\u00d78=\"],\n  [\"481\u00d72=\", \"351\u00d78=\"],\n  [\"128\u00d75=\", \"443\u00d77=\"],\n  [\"125\u00d78=\", \"724\u00d78=\"],\n  [\"771\u00d75=\", \"542\u00d75=\"],\n  [\"532\u00d77=\", \"799\u00d74=\"],\n  [\"835\u00d73=\", \"525\u00d78=\"],\n  [\"760\u00d75=\", \"758\u00d76=\"],\n  [\"142\u00d78=\", \"622\u00d78=\"],\n  [\"112\u00d78=\", \"349\u00d75=\"],\n  [\"770\u00d73=\", \"488\u00d79=\"],\n  [\"236\u00d73=\", \"444\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and all three-digit x one-digit multiplication\n# expressions in the practice sheet table to the new day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-07-10 Thursday\", \"2025-07-11 Friday\"),\n    @(\"122\u00d72=\", \"305\u00d73=\"),\n    @(\"868\u00d79=\", \"470\u00d79=\"),\n    @(\"327\u00d78=\", \"957\u00d77=\"),\n    @(\"200\u00d72=\", \"296\u00d73=\"),\n    @(\"116\u00d77=\", \"770\u00d78=\"),\n    @(\"685\u00d79=\", \"324\u00d79=\"),\n    @(\"891\u00d73=\", \"628\u00d75=\"),\n    @(\"603\u00d77=\", \"231\u00d76=\"),\n    @(\"964\u00d75=\", \"814\u00d78=\"),\n    @(\"322\u00d78=\", \"780\u00d75=\"),\n    @(\"681\u00d75=\", \"310\u00d72=\"),\n    @(\"587\u00d78=\", \"817\u00d79=\"),\n    @(\"325\u00d78=\", \"180\u00d74=\"),\n    @(\"602\u00d77=\", \"965\u00d78=\"),\n    @(\"481\u00d72=\", \"351\u00d78=\"),\n    @(\"128\u00d75=\", \"443\u00d77=\"),\n    @(\"125\u00d78=\", \"724\u00d78=\"),\n    @(\"771\u00d75=\", \"542\u00d75=\"),\n    @(\"532\u00d77=\", \"799\u00d74=\"),\n    @(\"835\u00d73=\", \"525\u00d78=\"),\n    @(\"760\u00d75=\", \"758\u00d76=\"),\n    @(\"142\u00d78=\", \"622\u00d78=\"),\n    @(\"112\u00d78=\", \"349\u00d75=\"),\n    @(\"770\u00d73=\", \"488\u00d79=\"),\n    @(\"236\u00d73=\", \"444\u00d76=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
